$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.53 = 17720.37 pesos`n✅ 17720.37 pesos = 4.52 = 974.78 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2.Range("N10").Value = 220.65
$ws2.Range("O10").Value = 3910
$ws2.Range("N12").Value = 3922.09
$ws2.Range("O12").Value = 215.75
